# "Add files via upload" -- the author re-uploaded Base_Cangaceiros_Empresas.xlsx
# with a refreshed "Empresas2" sheet: the company list grew from 15 to 41 rows
# (now alphabetised, with many new companies added) and the workbook was left
# open on the "Empresas2" tab instead of "Cangaceiros2".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Empresas2")

# Row 1 is the header; rows 2.. are the company data (columns A-E). A couple
# of rows only go out to column D (no "Localidade" recorded), so those inner
# arrays are simply shorter.
$data = @(
    @("Empresa", "Área de Atuação", "Maior Necessidade", "Segunda Necessidade", "Localidade"),
    @("59mil", "Comunicação e Criatividade", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Natal"),
    @("ACONT - ASSESSORIA E CONSULTORIA CONTÁBIL", "Engenharia e Tecnologia", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Natal"),
    @("Agro+ Assessoria e Consultoria Agrícola", "Engenharia e Tecnologia", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Natal"),
    @("ALPE Engenharia", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Mossoró"),
    @("ANIMUS Consultoria Jurídica", "Negócios e Gestão", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Natal"),
    @("Apex Empreendedorismo e Soluções Jurídicas", "Negócios e Gestão", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Mossoró"),
    @("Byte Seridó Júnior", "Engenharia e Tecnologia", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Mossoró"),
    @("Cápsula Júnior", "Engenharia e Tecnologia", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Natal"),
    @("Concrete", "Engenharia e Tecnologia", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Natal"),
    @("Conecta Solutions", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Mossoró"),
    @("CONSEJ - Consultoria Jurídica Júnior", "Negócios e Gestão", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Natal"),
    @("Contabilize Jr. Soluções Contábeis", "Negócios e Gestão", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Mossoró"),
    @("CORE Engenharia Biomedica", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Natal"),
    @("Econsul Consultoria Econômica", "Negócios e Gestão", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Natal"),
    @("Edifique Jr. - Arquitetura e Engenharia Civil", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Natal"),
    @("EJECT", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Natal"),
    @("Elysium Consultoria Odontológica", "Negócios e Gestão", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Natal"),
    @("EMBASA - Empresa de Biotecnologia Aplicada ao Semiárido", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Mossoró"),
    @("EMJUZ - EMPRESA JUNIOR DE ZOOTECNIA", "Negócios e Gestão", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Mossoró"),
    @("Enfasis Júnior", "Comunicação e Criatividade", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Natal"),
    @("FLOWLINE ENGENHARIA", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)"),
    @("Holos Consultoria Jr", "Negócios e Gestão", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Natal"),
    @("Honoris Consultoria Juridica Junior", "Negócios e Gestão", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Natal"),
    @("Include Engenharia", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Natal"),
    @("Lastro Consultoria e Investimentos", "Negócios e Gestão", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Mossoró"),
    @("LUMUS Engenharia", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Natal"),
    @("Mechanics Consultoria & Serviços", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Mossoró"),
    @("Nexum Consultoria Jurídica", "Negócios e Gestão", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Natal"),
    @("nuteq", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Natal"),
    @("Pilares Engenharia Júnior", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Caraúbas"),
    @("Pirâmides Topografia e Projetos", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Pau dos Ferros"),
    @("PROAQUA JR", "Engenharia e Tecnologia", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Mossoró"),
    @("Project Jr Consultoria", "Negócios e Gestão", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Mossoró"),
    @("Quatro Elementos", "Negócios e Gestão", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Mossoró"),
    @("Sinergy Júnior Consultoria e Projetos", "Engenharia e Tecnologia", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Mossoró"),
    @("Síntesis Jr - Consultoria e Projetos em Engenharia Química", "Engenharia e Tecnologia", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Mossoró"),
    @("Solidus Júnior", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Natal"),
    @("SOLIF - Engenharia de Energia", "Engenharia e Tecnologia", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Natal"),
    @("Spell JR - Assessoria em Língua Inglesa", "Comunicação e Criatividade", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Mossoró"),
    @("Unifica Floresta", "Engenharia e Tecnologia", "Time e Cultura (atração e retenção de membros, formação de lideranças, engajamento do time, cultura organizacional)", "Projetos e Modelo de Negócios (proposta de valor, execução de projetos, cadeia de valor, inovação em soluções)", "Mossoró"),
    @("Universitur", "Comunicação e Criatividade", "Gestão e Operações (planejamento estratégico, sistema de gestão, gestão financeira, estrutura organizacional)", "Vendas e Mercado (processo de vendas, estratégia comercial, aquisição de clientes, retenção de clientes)", "Natal")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}

# Columns C/D ("Maior Necessidade" / "Segunda Necessidade") move from a
# shared 108.71-char autofit width to two narrower, independent custom widths.
$ws.Columns.Item(3).ColumnWidth = 25
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666

# The workbook now opens on the "Empresas2" tab (previously "Cangaceiros2"),
# with the cursor left on E36.
$ws.Activate()
$ws.Range("E36").Select()
